$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(18, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2210", 31601)
    ,@(19, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2209", 35112)
    ,@(20, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2208", 35112)
    ,@(21, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2207", 35112)
    ,@(22, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2206", 35112)
    ,@(23, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2205", 35112)
    ,@(24, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2204", 35112)
    ,@(25, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2203", 35112)
    ,@(26, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2202", 35112)
    ,@(27, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2201", 35112)
    ,@(28, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2112", 35112)
    ,@(29, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2111", 35112)
    ,@(30, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2110", 35112)
    ,@(31, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2109", 35112)
    ,@(32, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2108", 35112)
    ,@(33, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2107", 35112)
    ,@(34, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2106", 35112)
    ,@(35, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2105", 35112)
    ,@(36, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2104", 35112)
    ,@(37, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2103", 35112)
    ,@(38, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2102", 35112)
    ,@(39, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2101", 35112)
    ,@(40, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2012", 35112)
    ,@(41, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2011", 35112)
    ,@(42, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2010", 35112)
    ,@(43, "1050965716", "JOSE ANGEL FLOREZ PUELLO", "2009", 35112)
    ,@(44, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2210", 31601)
    ,@(45, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2209", 35112)
    ,@(46, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2208", 35112)
    ,@(47, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2207", 35112)
    ,@(48, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2206", 35112)
    ,@(49, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2205", 35112)
    ,@(50, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2204", 35112)
    ,@(51, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2203", 35112)
    ,@(52, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2202", 35112)
    ,@(53, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2201", 35112)
    ,@(54, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2112", 35112)
    ,@(55, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2111", 35112)
    ,@(56, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2110", 35112)
    ,@(57, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2109", 35112)
    ,@(58, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2108", 35112)
    ,@(59, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2107", 35112)
    ,@(60, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2106", 35112)
    ,@(61, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2105", 35112)
    ,@(62, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2104", 35112)
    ,@(63, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2103", 35112)
    ,@(64, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2102", 35112)
    ,@(65, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2101", 35112)
    ,@(66, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2012", 35112)
    ,@(67, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2011", 35112)
    ,@(68, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2010", 35112)
    ,@(69, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2009", 35112)
    ,@(70, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2008", 35112)
    ,@(71, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2007", 35112)
    ,@(72, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2006", 35112)
    ,@(73, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2005", 35112)
    ,@(74, "1142949696", "DIANA CAROLINA INCIARTE MUÑOZ", "2004", 35112)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
}